$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# @@ -2775,25 +2775,25 @@
$ws.Range("H43").Value = 989.6
$ws.Range("I43").Value = 866.6667
$ws.Range("J43").Value = 1042.2858
$ws.Range("K43").Value = 866.6667
$ws.Range("L43").Value = 1042.2858
$ws.Range("M43").Value = -797.6667
$ws.Range("N43").Value = -1180.2858
# @@ -7265,25 +7265,25 @@
$ws.Range("H132").Value = 2028.3422
$ws.Range("I132").Value = 1574.0344
$ws.Range("J132").Value = 3492.2222
$ws.Range("K132").Value = 4722.1032
$ws.Range("L132").Value = 10476.6666
$ws.Range("M132").Value = -2192.1032
$ws.Range("N132").Value = -15536.6666
# @@ -7568,25 +7568,25 @@
$ws.Range("H138").Value = 2130973
$ws.Range("J138").Value = 3732.4
$ws.Range("L138").Value = 11197.2
$ws.Range("N138").Value = -21477.2

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# @@ -7867,25 +7867,25 @@
$ws.Range("H2").Value = 1156.2972
$ws.Range("I2").Value = 716.4400000000001
$ws.Range("J2").Value = 2072.6667
$ws.Range("K2").Value = 716.4400000000001
$ws.Range("L2").Value = 2072.6667
$ws.Range("M2").Value = -603.4400000000001
$ws.Range("N2").Value = -2298.6667
# @@ -11425,22 +11425,22 @@
$ws.Range("H74").Value = 1083.6875
$ws.Range("I74").Value = 981.48
$ws.Range("K74").Value = 981.48
$ws.Range("M74").Value = -107.48
# @@ -11572,22 +11572,22 @@
$ws.Range("H77").Value = 1083.6875
$ws.Range("I77").Value = 981.48
$ws.Range("K77").Value = 4907.4
$ws.Range("M77").Value = -539.3999999999996
# @@ -12999,22 +12999,22 @@
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524
# @@ -13483,25 +13483,25 @@
$ws.Range("H116").Value = 1156.2972
$ws.Range("I116").Value = 716.4400000000001
$ws.Range("J116").Value = 2072.6667
$ws.Range("K116").Value = 716.4400000000001
$ws.Range("L116").Value = 2072.6667
$ws.Range("M116").Value = 1577.56
$ws.Range("N116").Value = -6660.6667
# @@ -14258,25 +14258,25 @@
$ws.Range("H132").Value = 2477.2273
$ws.Range("I132").Value = 1998.9166
$ws.Range("J132").Value = 3051.2
$ws.Range("K132").Value = 5996.7498
$ws.Range("L132").Value = 9153.599999999999
$ws.Range("M132").Value = -3466.7498
$ws.Range("N132").Value = -14213.6

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# @@ -14894,25 +14894,25 @@
$ws.Range("H3").Value = 1156.2972
$ws.Range("I3").Value = 716.4400000000001
$ws.Range("J3").Value = 2072.6667
$ws.Range("K3").Value = 716.4400000000001
$ws.Range("L3").Value = 2072.6667
$ws.Range("M3").Value = -602.4400000000001
$ws.Range("N3").Value = -2300.6667
# @@ -19362,25 +19362,25 @@
$ws.Range("H94").Value = 1883.3334
$ws.Range("I94").Value = 1133.3334
$ws.Range("J94").Value = 2633.3333
$ws.Range("K94").Value = 1133.3334
$ws.Range("L94").Value = 2633.3333
$ws.Range("M94").Value = -682.3334
$ws.Range("N94").Value = -3535.3333
# @@ -21313,25 +21313,25 @@
$ws.Range("H134").Value = 2782.15
$ws.Range("I134").Value = 2437.1667
$ws.Range("J134").Value = 3299.625
$ws.Range("K134").Value = 7311.500100000001
$ws.Range("L134").Value = 9898.875
$ws.Range("M134").Value = -4776.500100000001
$ws.Range("N134").Value = -14968.875

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# @@ -21900,25 +21900,22 @@
$ws.Range("H4").Value = 9833.333000000001
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 9833.333000000001
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 9833.333000000001
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -10057.333
# @@ -26361,25 +26358,25 @@
$ws.Range("H94").Value = 10970.765
$ws.Range("I94").Value = 1172
$ws.Range("J94").Value = 13985.77
$ws.Range("K94").Value = 1172
$ws.Range("L94").Value = 13985.77
$ws.Range("M94").Value = -721
$ws.Range("N94").Value = -14887.77
# @@ -26995,25 +26992,25 @@
$ws.Range("H107").Value = 440.66666
$ws.Range("I107").Value = 426.33334
$ws.Range("J107").Value = 455
$ws.Range("K107").Value = 426.33334
$ws.Range("L107").Value = 455
$ws.Range("M107").Value = 1493.66666
$ws.Range("N107").Value = -4295
# @@ -28226,25 +28223,25 @@
$ws.Range("H132").Value = 1849.8572
$ws.Range("I132").Value = 1381.1875
$ws.Range("J132").Value = 3349.6
$ws.Range("K132").Value = 4143.5625
$ws.Range("L132").Value = 10048.8
$ws.Range("M132").Value = -1613.5625
$ws.Range("N132").Value = -15108.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# @@ -29887,25 +29884,25 @@
$ws.Range("H23").Value = 455
$ws.Range("J23").Value = 440
$ws.Range("L23").Value = 1320
$ws.Range("N23").Value = -1790
# @@ -35426,25 +35423,25 @@
$ws.Range("H132").Value = 1506.2609
$ws.Range("I132").Value = 1008.7143
$ws.Range("J132").Value = 1723.9375
$ws.Range("K132").Value = 9078.4287
$ws.Range("L132").Value = 15515.4375
$ws.Range("M132").Value = -6548.4287
$ws.Range("N132").Value = -20575.4375

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# @@ -42407,25 +42404,25 @@
$ws.Range("H132").Value = 1992.5
$ws.Range("I132").Value = 1482.4166
$ws.Range("J132").Value = 4032.8333
$ws.Range("K132").Value = 4447.2498
$ws.Range("L132").Value = 12098.4999
$ws.Range("M132").Value = -1917.2498
$ws.Range("N132").Value = -17158.4999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# @@ -45144,25 +45141,25 @@
$ws.Range("H46").Value = 1559.2
$ws.Range("I46").Value = 1566.6666
$ws.Range("J46").Value = 1548
$ws.Range("K46").Value = 1566.6666
$ws.Range("L46").Value = 1548
$ws.Range("M46").Value = -1378.6666
$ws.Range("N46").Value = -1924
# @@ -45573,25 +45570,25 @@
$ws.Range("H55").Value = 136
$ws.Range("I55").Value = 167
$ws.Range("J55").Value = 89.5
$ws.Range("K55").Value = 167
$ws.Range("L55").Value = 89.5
$ws.Range("M55").Value = 6
$ws.Range("N55").Value = -435.5
# @@ -45870,25 +45867,25 @@
$ws.Range("H61").Value = 17130.285
$ws.Range("I61").Value = 27964
$ws.Range("J61").Value = 2685.3333
$ws.Range("K61").Value = 27964
$ws.Range("L61").Value = 2685.3333
$ws.Range("M61").Value = -27762
$ws.Range("N61").Value = -3089.3333
# @@ -48412,25 +48409,25 @@
$ws.Range("H113").Value = 17130.285
$ws.Range("I113").Value = 27964
$ws.Range("J113").Value = 2685.3333
$ws.Range("K113").Value = 27964
$ws.Range("L113").Value = 2685.3333
$ws.Range("M113").Value = -25794
$ws.Range("N113").Value = -7025.3333
# @@ -49343,25 +49340,25 @@
$ws.Range("H132").Value = 5175.375
$ws.Range("I132").Value = 5265.5884
$ws.Range("J132").Value = 4956.2856
$ws.Range("K132").Value = 15796.7652
$ws.Range("L132").Value = 14868.8568
$ws.Range("M132").Value = -13266.7652
$ws.Range("N132").Value = -19928.8568
# @@ -49542,22 +49539,22 @@
$ws.Range("H136").Value = 2917.1365
$ws.Range("I136").Value = 2508.85
$ws.Range("K136").Value = 7526.549999999999
$ws.Range("M136").Value = -4976.549999999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# @@ -53892,19 +53889,22 @@
$ws.Range("H82").Value = 30150.5
$ws.Range("J82").Value = 30150.5
$ws.Range("L82").Value = 30150.5
$ws.Range("N82").Value = -30916.5
# @@ -54036,19 +54036,22 @@
$ws.Range("H85").Value = 30150.5
$ws.Range("J85").Value = 30150.5
$ws.Range("L85").Value = 30150.5
$ws.Range("N85").Value = -32802.5
# @@ -56348,25 +56351,25 @@
$ws.Range("H132").Value = 1437.4073
$ws.Range("I132").Value = 1090.6
$ws.Range("J132").Value = 2428.2856
$ws.Range("K132").Value = 3271.8
$ws.Range("L132").Value = 7284.8568
$ws.Range("M132").Value = -741.7999999999997
$ws.Range("N132").Value = -12344.8568
# @@ -56648,22 +56651,22 @@
$ws.Range("H138").Value = 78997.5
$ws.Range("J138").Value = 78997.5
$ws.Range("L138").Value = 78997.5
$ws.Range("N138").Value = -89277.5
